$p = $ppt.ActivePresentation
$s = $p.Slides.Item(7)
$sh = $s.Shapes.Item(3)
$tf = $sh.TextFrame
$tr = $tf.TextRange

# Preserve the shape's current (auto-fit) height: the upcoming text-content
# tweak would otherwise make PowerPoint re-flow/re-measure this text box.
$origHeight = $sh.Height

# Paragraph 3 of this textbox currently reads:
#   "인라인 스타일은 가장 우선순위가 제일 큰 스타일이다.(내부스타일 <인라인스타일)"
# Split the leading sentence into four runs so that two of the pieces
# ("인라인 스타일은 " and "제일 큰 스타일") become red, while
# "가장 우선순위가 " and "이다" keep the default color.
$para3 = $tr.Paragraphs(3, 1)

# "인라인 스타일은 " -> red (FF0000)
$para3.Characters(1, 9).Font.Color.RGB = 255
# "제일 큰 스타일" -> red (FF0000)
$para3.Characters(19, 8).Font.Color.RGB = 255

# Drop the trailing space after "내부스타일" later in the same paragraph
# (was "내부스타일 ", now "내부스타일").
$para3.Characters(31, 6).Text = "내부스타일"

# Restore the text box's height so the auto-fit re-layout triggered by the
# text edit above doesn't leave a stray size change behind.
$sh.Height = $origHeight
